$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: fix the first three summary cells -----------------------------
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# --- Step 2: insert the 10 new per-iteration rows right after row 3 --------
$newValues = @("368","0.00002","0.00005","0.00002","0.00001","0.00003","0.00003","0.00004","0.00990","100.0")

$insertBefore = $t.Rows.Item(4)
foreach ($val in $newValues) {
    $t.Rows.Add($insertBefore) | Out-Null
}

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $t.Cell(4 + $i, 1).Range.Text = $newValues[$i]
}

# --- Step 3: collapse the three tab-separated summary rows back down to a
#             single value each (they now sit after the rows we just added)
$rowCount = $t.Rows.Count
$t.Cell($rowCount - 2,1).Range.Text = "99.99"
$t.Cell($rowCount - 1,1).Range.Text = "0.01"
$t.Cell($rowCount,1).Range.Text = "153"
